# Update "想去人数" (want-to-go count) figures in the F column
# on both the "展览" (Exhibition) sheet and the "全部类型" (All types) sheet.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 5561
$wsExhibition.Range("F3").Value = 170
$wsExhibition.Range("F4").Value = 948
$wsExhibition.Range("F5").Value = 12

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 5561
$wsAllTypes.Range("F3").Value = 170
$wsAllTypes.Range("F4").Value = 948
$wsAllTypes.Range("F5").Value = 12
